$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1470.3469
$ws.Cells.Item(17, 10).Value = 1575.5116
$ws.Cells.Item(17, 12).Value = 4726.5348
$ws.Cells.Item(17, 14).Value = -5062.5348

$ws.Cells.Item(40, 8).Value = 5254.8125
$ws.Cells.Item(40, 9).Value = 3871.1428
$ws.Cells.Item(40, 10).Value = 6331
$ws.Cells.Item(40, 11).Value = 3871.1428
$ws.Cells.Item(40, 12).Value = 6331
$ws.Cells.Item(40, 13).Value = -3696.1428
$ws.Cells.Item(40, 14).Value = -6681

$ws.Cells.Item(43, 8).Value = 221271.05
$ws.Cells.Item(43, 10).Value = 320384.62
$ws.Cells.Item(43, 12).Value = 320384.62
$ws.Cells.Item(43, 14).Value = -320522.62

$ws.Cells.Item(86, 8).Value = 81849460
$ws.Cells.Item(86, 9).Value = 82409200
$ws.Cells.Item(86, 11).Value = 82409200
$ws.Cells.Item(86, 13).Value = -82408077

$ws.Cells.Item(89, 8).Value = 81849460
$ws.Cells.Item(89, 9).Value = 82409200
$ws.Cells.Item(89, 11).Value = 412046000
$ws.Cells.Item(89, 13).Value = -412040384

$ws.Cells.Item(107, 8).Value = 31252052
$ws.Cells.Item(107, 9).Value = 12502422
$ws.Cells.Item(107, 10).Value = 125000200
$ws.Cells.Item(107, 11).Value = 12502422
$ws.Cells.Item(107, 12).Value = 125000200
$ws.Cells.Item(107, 13).Value = -12500502
$ws.Cells.Item(107, 14).Value = -125004040

$ws.Cells.Item(112, 8).Value = 2797.705
$ws.Cells.Item(112, 10).Value = 2797.705
$ws.Cells.Item(112, 12).Value = 8393.115
$ws.Cells.Item(112, 14).Value = -10609.115

$ws.Cells.Item(129, 8).Value = 1145.5625
$ws.Cells.Item(129, 9).Value = 693
$ws.Cells.Item(129, 10).Value = 1899.8334
$ws.Cells.Item(129, 11).Value = 2079
$ws.Cells.Item(129, 12).Value = 5699.5002
$ws.Cells.Item(129, 13).Value = 2921
$ws.Cells.Item(129, 14).Value = -15699.5002

$ws.Cells.Item(137, 8).Value = 4046.2666
$ws.Cells.Item(137, 9).Value = 4437.75
$ws.Cells.Item(137, 10).Value = 3903.9092
$ws.Cells.Item(137, 11).Value = 13313.25
$ws.Cells.Item(137, 12).Value = 11711.7276
$ws.Cells.Item(137, 13).Value = -10763.25
$ws.Cells.Item(137, 14).Value = -16811.7276

$ws.Cells.Item(138, 8).Value = 1474110.1
$ws.Cells.Item(138, 9).Value = 997.32355
$ws.Cells.Item(138, 10).Value = 2947223
$ws.Cells.Item(138, 11).Value = 2991.97065
$ws.Cells.Item(138, 12).Value = 8841669
$ws.Cells.Item(138, 13).Value = 2148.02935
$ws.Cells.Item(138, 14).Value = -8851949

$ws.Cells.Item(141, 8).Value = 1393.3914
$ws.Cells.Item(141, 9).Value = 1229.6364
$ws.Cells.Item(141, 10).Value = 4996
$ws.Cells.Item(141, 11).Value = 3688.9092
$ws.Cells.Item(141, 12).Value = 14988
$ws.Cells.Item(141, 13).Value = 1491.0908
$ws.Cells.Item(141, 14).Value = -25348

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1871979
$ws.Cells.Item(32, 9).Value = 2087796.4
$ws.Cells.Item(32, 11).Value = 2087796.4
$ws.Cells.Item(32, 13).Value = -2087509.4

$ws.Cells.Item(61, 8).Value = 5734.436
$ws.Cells.Item(61, 9).Value = 2173.6296
$ws.Cells.Item(61, 10).Value = 13746.25
$ws.Cells.Item(61, 11).Value = 2173.6296
$ws.Cells.Item(61, 12).Value = 13746.25
$ws.Cells.Item(61, 13).Value = -1961.6296
$ws.Cells.Item(61, 14).Value = -14170.25

$ws.Cells.Item(74, 8).Value = 30025.174
$ws.Cells.Item(74, 9).Value = 43385.332
$ws.Cells.Item(74, 11).Value = 43385.332
$ws.Cells.Item(74, 13).Value = -42511.332

$ws.Cells.Item(77, 8).Value = 30025.174
$ws.Cells.Item(77, 9).Value = 43385.332
$ws.Cells.Item(77, 11).Value = 216926.66
$ws.Cells.Item(77, 13).Value = -212558.66

$ws.Cells.Item(103, 8).Value = 54657.332
$ws.Cells.Item(103, 10).Value = 54657.332
$ws.Cells.Item(103, 12).Value = 54657.332
$ws.Cells.Item(103, 14).Value = -57001.332

$ws.Cells.Item(136, 8).Value = 5734.436
$ws.Cells.Item(136, 9).Value = 2173.6296
$ws.Cells.Item(136, 10).Value = 13746.25
$ws.Cells.Item(136, 11).Value = 6520.888800000001
$ws.Cells.Item(136, 12).Value = 41238.75
$ws.Cells.Item(136, 13).Value = -3970.888800000001
$ws.Cells.Item(136, 14).Value = -46338.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 248
$ws.Cells.Item(22, 9).Value = 248
$ws.Cells.Item(22, 11).Value = 248
$ws.Cells.Item(22, 13).Value = -75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7100051
$ws.Cells.Item(31, 9).Value = 4087.1667
$ws.Cells.Item(31, 10).Value = 11504443
$ws.Cells.Item(31, 11).Value = 4087.1667
$ws.Cells.Item(31, 12).Value = 11504443
$ws.Cells.Item(31, 13).Value = -3792.1667
$ws.Cells.Item(31, 14).Value = -11505033

$ws.Cells.Item(34, 8).Value = 7100051
$ws.Cells.Item(34, 9).Value = 4087.1667
$ws.Cells.Item(34, 10).Value = 11504443
$ws.Cells.Item(34, 11).Value = 4087.1667
$ws.Cells.Item(34, 12).Value = 11504443
$ws.Cells.Item(34, 13).Value = -3885.1667
$ws.Cells.Item(34, 14).Value = -11504847

$ws.Cells.Item(58, 8).Value = 10422838
$ws.Cells.Item(58, 9).Value = 21741896
$ws.Cells.Item(58, 11).Value = 21741896
$ws.Cells.Item(58, 13).Value = -21741693

$ws.Cells.Item(107, 8).Value = 2203.1853
$ws.Cells.Item(107, 9).Value = 2209.7
$ws.Cells.Item(107, 11).Value = 2209.7
$ws.Cells.Item(107, 13).Value = -289.6999999999998

$ws.Cells.Item(132, 8).Value = 7848954
$ws.Cells.Item(132, 9).Value = 3718.5833
$ws.Cells.Item(132, 11).Value = 11155.7499
$ws.Cells.Item(132, 13).Value = -8625.749899999999

$ws.Cells.Item(134, 8).Value = 5770.58
$ws.Cells.Item(134, 9).Value = 2574.2173
$ws.Cells.Item(134, 11).Value = 7722.651899999999
$ws.Cells.Item(134, 13).Value = -5187.651899999999

$ws.Cells.Item(136, 8).Value = 10422838
$ws.Cells.Item(136, 9).Value = 21741896
$ws.Cells.Item(136, 11).Value = 65225688
$ws.Cells.Item(136, 13).Value = -65223138

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(112, 8).Value = 3547.25
$ws.Cells.Item(112, 9).Value = 2796.3333
$ws.Cells.Item(112, 11).Value = 8388.999899999999
$ws.Cells.Item(112, 13).Value = -7280.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 19333.334
$ws.Cells.Item(33, 10).Value = 19333.334
$ws.Cells.Item(33, 12).Value = 19333.334
$ws.Cells.Item(33, 14).Value = -19837.334

$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 14).ClearContents()

$ws.Cells.Item(40, 8).Value = 22000

$ws.Cells.Item(43, 8).Value = 3000
$ws.Cells.Item(43, 9).Value = 3000
$ws.Cells.Item(43, 11).Value = 3000
$ws.Cells.Item(43, 13).Value = -2849

$ws.Cells.Item(97, 8).Value = 1538.0588
$ws.Cells.Item(97, 9).Value = 1166.091
$ws.Cells.Item(97, 10).Value = 2220
$ws.Cells.Item(97, 11).Value = 1166.091
$ws.Cells.Item(97, 12).Value = 2220
$ws.Cells.Item(97, 13).Value = -670.0909999999999
$ws.Cells.Item(97, 14).Value = -3212

$ws.Cells.Item(102, 8).Value = 1426.3771
$ws.Cells.Item(102, 9).Value = 1240.2037
$ws.Cells.Item(102, 11).Value = 1240.2037
$ws.Cells.Item(102, 13).Value = 381.7963

$ws.Cells.Item(132, 8).Value = 4763.702
$ws.Cells.Item(132, 9).Value = 2447.476
$ws.Cells.Item(132, 10).Value = 6634.5
$ws.Cells.Item(132, 11).Value = 7342.428
$ws.Cells.Item(132, 12).Value = 19903.5
$ws.Cells.Item(132, 13).Value = -4812.428
$ws.Cells.Item(132, 14).Value = -24963.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3263.4375
$ws.Cells.Item(22, 9).Value = 1980
$ws.Cells.Item(22, 11).Value = 1980
$ws.Cells.Item(22, 13).Value = -1685

$ws.Cells.Item(27, 8).Value = 3263.4375
$ws.Cells.Item(27, 9).Value = 1980
$ws.Cells.Item(27, 11).Value = 1980
$ws.Cells.Item(27, 13).Value = -1873

$ws.Cells.Item(40, 8).Value = 5500.227
$ws.Cells.Item(40, 9).Value = 3218.2727
$ws.Cells.Item(40, 11).Value = 3218.2727
$ws.Cells.Item(40, 13).Value = -3082.2727

$ws.Cells.Item(46, 8).Value = 6947716
$ws.Cells.Item(46, 9).Value = 1733.3334
$ws.Cells.Item(46, 10).Value = 11115306
$ws.Cells.Item(46, 11).Value = 1733.3334
$ws.Cells.Item(46, 12).Value = 11115306
$ws.Cells.Item(46, 13).Value = -1545.3334
$ws.Cells.Item(46, 14).Value = -11115682

$ws.Cells.Item(122, 8).Value = 5786.174
$ws.Cells.Item(122, 9).Value = 2866.889
$ws.Cells.Item(122, 11).Value = 8600.667000000001
$ws.Cells.Item(122, 13).Value = -6150.667000000001

$ws.Cells.Item(132, 8).Value = 8626760
$ws.Cells.Item(132, 9).Value = 17859830
$ws.Cells.Item(132, 10).Value = 9228
$ws.Cells.Item(132, 11).Value = 53579490
$ws.Cells.Item(132, 12).Value = 27684
$ws.Cells.Item(132, 13).Value = -53576960
$ws.Cells.Item(132, 14).Value = -32744

$ws.Cells.Item(136, 8).Value = 5895.107
$ws.Cells.Item(136, 9).Value = 1823.0605
$ws.Cells.Item(136, 11).Value = 5469.181500000001
$ws.Cells.Item(136, 13).Value = -2919.181500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 102573.4
$ws.Cells.Item(122, 9).Value = 139796.17
$ws.Cells.Item(122, 10).Value = 4440.636
$ws.Cells.Item(122, 11).Value = 419388.51
$ws.Cells.Item(122, 12).Value = 13321.908
$ws.Cells.Item(122, 13).Value = -416938.51
$ws.Cells.Item(122, 14).Value = -18221.908

$ws.Cells.Item(136, 8).Value = 19630998
$ws.Cells.Item(136, 9).Value = 34483496
$ws.Cells.Item(136, 11).Value = 103450488
$ws.Cells.Item(136, 13).Value = -103447938
